$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'274.47"
$ws.Range("E2").Value = "'-1.36%"

$ws.Range("D3").Value = "'26.77"
$ws.Range("E3").Value = "'-2.27%"

$ws.Range("D4").Value = "'4.849"
$ws.Range("E4").Value = "'0.16%"

$ws.Range("D5").Value = "'0.06323"
$ws.Range("E5").Value = "'1.17%"

$ws.Range("D6").Value = "'6.884"
$ws.Range("E6").Value = "'-0.23%"

$ws.Range("D7").Value = "'3.322"
$ws.Range("E7").Value = "'1.69%"

$ws.Range("D8").Value = "'1.280"
$ws.Range("E8").Value = "'35.33%"

$ws.Range("D9").Value = "'0.8707"
$ws.Range("E9").Value = "'-0.90%"

$ws.Range("E10").Value = "'0.80%"

$ws.Range("D11").Value = "'0.04994"
$ws.Range("E11").Value = "'-4.37%"

$ws.Range("D12").Value = "'0.07402"
$ws.Range("E12").Value = "'1.50%"

$ws.Range("D13").Value = "'0.02941"
$ws.Range("E13").Value = "'-6.92%"

$ws.Range("D14").Value = "'0.09027"
$ws.Range("E14").Value = "'-0.32%"

$ws.Range("D15").Value = "'0.001569"
$ws.Range("E15").Value = "'1.09%"

$ws.Range("D16").Value = "'0.0006315"
$ws.Range("E16").Value = "'0.76%"

$ws.Range("D17").Value = "'0.006020"
$ws.Range("E17").Value = "'-1.29%"

$ws.Range("D18").Value = "'3.445"
$ws.Range("E18").Value = "'-0.42%"

$ws.Range("D19").Value = "'2.283"
$ws.Range("E19").Value = "'-0.08%"

$ws.Range("E20").Value = "'1.28%"

$ws.Range("D21").Value = "'0.1335"
$ws.Range("E21").Value = "'1.89%"

$ws.Range("D22").Value = "'3.916"
$ws.Range("E22").Value = "'1.72%"

$ws.Range("D23").Value = "'0.04348"
$ws.Range("E23").Value = "'0.50%"

$ws.Range("D24").Value = "'0.001179"
$ws.Range("E24").Value = "'0.36%"

$ws.Range("D25").Value = "'0.004254"
$ws.Range("E25").Value = "'-0.49%"

$ws.Range("E26").Value = "'0.09%"

$ws.Range("D40").Value = "'0.04045"
$ws.Range("E40").Value = "'0.44%"

$ws.Range("D41").Value = "'0.006687"
$ws.Range("E41").Value = "'4.30%"

$ws.Range("D42").Value = "'0.1164"
$ws.Range("E42").Value = "'0.90%"

$ws.Range("D43").Value = "'0.002100"
$ws.Range("E43").Value = "'-0.42%"

$ws.Range("E44").Value = "'-10.79%"

$ws.Range("D45").Value = "'0.00005304"
$ws.Range("E45").Value = "'4.26%"

$ws.Range("E46").Value = "'-37.42%"
